$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5833103656768799
$ws.Range("B1").Value = 3.550758600234985
$ws.Range("C1").Value = 4.359034061431885
$ws.Range("D1").Value = 2.565619468688965
$ws.Range("E1").Value = 1.060403227806091
